$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: assign a value while forcing it to stay a text cell (so
# numeric-looking strings like "1.00" or "11.98" keep their exact
# formatting/trailing zeros instead of Excel coercing them into a
# number). The temporary "@" number format is reverted to Normal
# immediately after so the cell keeps its original (unstyled) look.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "34.731.23"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "1.789.47"
$ws.Range("E3").Value = "  +0.70%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "223.01"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  -0.03%  "
Set-TextValue $ws.Range("D8") "32.72"
$ws.Range("E8").Value = "  +7.77%  "
$ws.Range("E9").Value = "  +0.86%  "
Set-TextValue $ws.Range("D10") "0.0686"
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "2.046.85"
$ws.Range("E12").Value = "  +0.75%  "
Set-TextValue $ws.Range("D13") "11.06"
$ws.Range("D14").Value = "1.790.39"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "34.750.53"
$ws.Range("E15").Value = "  +3.01%  "
Set-TextValue $ws.Range("D16") "0.631"
$ws.Range("E16").Value = "  +1.21%  "
Set-TextValue $ws.Range("D17") "4.31"
$ws.Range("E17").Value = "  +3.33%  "
Set-TextValue $ws.Range("D18") "68.49"
$ws.Range("E18").Value = "  +0.03%  "
Set-TextValue $ws.Range("D19") "253.31"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").Value = "0.0₃0783"
$ws.Range("E20").Value = "  +6.14%  "
Set-TextValue $ws.Range("D22") "10.48"
$ws.Range("E22").Value = "  +2.20%  "
Set-TextValue $ws.Range("D23") "4.19"
$ws.Range("E24").Value = "  -0.98%  "
Set-TextValue $ws.Range("D25") "158.82"
$ws.Range("E25").Value = "  +0.35%  "
Set-TextValue $ws.Range("D26") "16.35"
$ws.Range("E26").Value = "  -0.39%  "
Set-TextValue $ws.Range("D27") "7.06"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D30") "0.0515"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D31") "3.75"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Value = "1.437.23"
$ws.Range("E35").Value = "  -2.81%  "
Set-TextValue $ws.Range("D36") "1.06"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "0.632"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.0189"
$ws.Range("E38").Value = "  +2.52%  "
Set-TextValue $ws.Range("D39") "82.87"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("E41").Value = "  +0.07%  "
Set-TextValue $ws.Range("D42") "0.904"
Set-TextValue $ws.Range("D43") "2.06"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.07"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D45") "0.0503"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("E46").Value = "  +4.46%  "
$ws.Range("D47").Value = "1.943.78"
$ws.Range("E47").Value = "  +0.75%  "
Set-TextValue $ws.Range("D48") "104.82"
$ws.Range("E48").Value = "  +7.74%  "
Set-TextValue $ws.Range("D49") "11.98"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("E50").Value = "  -0.09%  "
Set-TextValue $ws.Range("D51") "49.63"
$ws.Range("E51").Value = "  -2.43%  "
